$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(18, 8).Value = 1422
$ws.Cells.Item(18, 9).Value = 1422
$ws.Cells.Item(18, 11).Value = 1422
$ws.Cells.Item(18, 13).Value = -1138
$ws.Cells.Item(51, 8).Value = 2965.5
$ws.Cells.Item(51, 9).Value = 2980.303
$ws.Cells.Item(51, 10).Value = 1500
$ws.Cells.Item(51, 11).Value = 2980.303
$ws.Cells.Item(51, 12).Value = 1500
$ws.Cells.Item(51, 13).Value = -2496.303
$ws.Cells.Item(51, 14).Value = -2468
$ws.Cells.Item(87, 8).Value = 66000
$ws.Cells.Item(87, 10).Value = 66000
$ws.Cells.Item(87, 12).Value = 66000
$ws.Cells.Item(87, 14).Value = -68496
$ws.Cells.Item(90, 8).Value = 66000
$ws.Cells.Item(90, 10).Value = 66000
$ws.Cells.Item(90, 12).Value = 198000
$ws.Cells.Item(90, 14).Value = -210480
$ws.Cells.Item(112, 8).Value = 1852.5294
$ws.Cells.Item(112, 9).Value = 1998.75
$ws.Cells.Item(112, 10).Value = 1807.5385
$ws.Cells.Item(112, 11).Value = 5996.25
$ws.Cells.Item(112, 12).Value = 5422.6155
$ws.Cells.Item(112, 13).Value = -4888.25
$ws.Cells.Item(112, 14).Value = -7638.6155
$ws.Cells.Item(138, 8).Value = 6857.4287
$ws.Cells.Item(138, 10).Value = 6333.6665
$ws.Cells.Item(138, 12).Value = 19000.9995
$ws.Cells.Item(138, 14).Value = -29280.9995
$ws.Cells.Item(141, 8).Value = 12228.974
$ws.Cells.Item(141, 9).Value = 4470.1
$ws.Cells.Item(141, 11).Value = 13410.3
$ws.Cells.Item(141, 13).Value = -8230.300000000001

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(17, 8).Value = 1750
$ws.Cells.Item(17, 9).Value = 1800
$ws.Cells.Item(17, 10).Value = 1700
$ws.Cells.Item(17, 11).Value = 1800
$ws.Cells.Item(17, 12).Value = 1700
$ws.Cells.Item(17, 13).Value = -1627
$ws.Cells.Item(17, 14).Value = -2046
$ws.Cells.Item(61, 8).Value = 2227518.2
$ws.Cells.Item(61, 9).Value = 5746.364
$ws.Cells.Item(61, 10).Value = 8337391
$ws.Cells.Item(61, 11).Value = 5746.364
$ws.Cells.Item(61, 12).Value = 8337391
$ws.Cells.Item(61, 13).Value = -5534.364
$ws.Cells.Item(61, 14).Value = -8337815
$ws.Cells.Item(74, 8).Value = 1695444.9
$ws.Cells.Item(74, 9).Value = 1925464.9
$ws.Cells.Item(74, 11).Value = 1925464.9
$ws.Cells.Item(74, 13).Value = -1924590.9
$ws.Cells.Item(77, 8).Value = 1695444.9
$ws.Cells.Item(77, 9).Value = 1925464.9
$ws.Cells.Item(77, 11).Value = 9627324.5
$ws.Cells.Item(77, 13).Value = -9622956.5
$ws.Cells.Item(136, 8).Value = 2227518.2
$ws.Cells.Item(136, 9).Value = 5746.364
$ws.Cells.Item(136, 10).Value = 8337391
$ws.Cells.Item(136, 11).Value = 17239.092
$ws.Cells.Item(136, 12).Value = 25012173
$ws.Cells.Item(136, 13).Value = -14689.092
$ws.Cells.Item(136, 14).Value = -25017273

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(35, 8).Value = 31200
$ws.Cells.Item(35, 10).Value = 31200
$ws.Cells.Item(35, 12).Value = 31200
$ws.Cells.Item(35, 14).Value = -31820
$ws.Cells.Item(105, 8).Value = 3475.6
$ws.Cells.Item(105, 10).Value = 5794.875
$ws.Cells.Item(105, 12).Value = 5794.875
$ws.Cells.Item(105, 14).Value = -9288.875
$ws.Cells.Item(134, 8).Value = 4172140.5
$ws.Cells.Item(134, 9).Value = 6823.1333
$ws.Cells.Item(134, 10).Value = 16668092
$ws.Cells.Item(134, 11).Value = 20469.3999
$ws.Cells.Item(134, 12).Value = 50004276
$ws.Cells.Item(134, 13).Value = -17934.3999
$ws.Cells.Item(134, 14).Value = -50009346

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 112.052635
$ws.Cells.Item(7, 10).Value = 197.33333
$ws.Cells.Item(7, 12).Value = 197.33333
$ws.Cells.Item(7, 14).Value = -423.33333
$ws.Cells.Item(22, 8).Value = 2390.7144
$ws.Cells.Item(22, 9).Value = 907.8570999999999
$ws.Cells.Item(22, 11).Value = 907.8570999999999
$ws.Cells.Item(22, 13).Value = -557.8570999999999
$ws.Cells.Item(41, 8).Value = 23599.4
$ws.Cells.Item(41, 10).Value = 23599.4
$ws.Cells.Item(41, 12).Value = 23599.4
$ws.Cells.Item(41, 14).Value = -24455.4
$ws.Cells.Item(50, 8).Value = 29499.666
$ws.Cells.Item(50, 10).Value = 29499.666
$ws.Cells.Item(50, 12).Value = 29499.666
$ws.Cells.Item(50, 14).Value = -30749.666
$ws.Cells.Item(51, 8).Value = 29199.666
$ws.Cells.Item(51, 10).Value = 29199.666
$ws.Cells.Item(51, 12).Value = 29199.666
$ws.Cells.Item(51, 14).Value = -30671.666
$ws.Cells.Item(59, 8).Value = 37018.168
$ws.Cells.Item(59, 10).Value = 40528.5
$ws.Cells.Item(59, 12).Value = 40528.5
$ws.Cells.Item(59, 14).Value = -42818.5
$ws.Cells.Item(60, 8).Value = 13628.429
$ws.Cells.Item(60, 10).Value = 14399.833
$ws.Cells.Item(60, 12).Value = 14399.833
$ws.Cells.Item(60, 14).Value = -15421.833
$ws.Cells.Item(61, 8).Value = 29199.666
$ws.Cells.Item(61, 10).Value = 29199.666
$ws.Cells.Item(61, 12).Value = 29199.666
$ws.Cells.Item(61, 14).Value = -29895.666
$ws.Cells.Item(68, 8).Value = 40332.777
$ws.Cells.Item(68, 10).Value = 40332.777
$ws.Cells.Item(68, 12).Value = 40332.777
$ws.Cells.Item(68, 14).Value = -41830.777
$ws.Cells.Item(71, 8).Value = 40332.777
$ws.Cells.Item(71, 10).Value = 40332.777
$ws.Cells.Item(71, 12).Value = 120998.331
$ws.Cells.Item(71, 14).Value = -128486.331
$ws.Cells.Item(74, 8).Value = 33875.5
$ws.Cells.Item(74, 10).Value = 33875.5
$ws.Cells.Item(74, 12).Value = 33875.5
$ws.Cells.Item(74, 14).Value = -35623.5
$ws.Cells.Item(77, 8).Value = 33875.5
$ws.Cells.Item(77, 10).Value = 33875.5
$ws.Cells.Item(77, 12).Value = 101626.5
$ws.Cells.Item(77, 14).Value = -110362.5
$ws.Cells.Item(86, 8).Value = 31449.611
$ws.Cells.Item(86, 9).Value = 57920.734
$ws.Cells.Item(86, 11).Value = 57920.734
$ws.Cells.Item(86, 13).Value = -56797.734
$ws.Cells.Item(89, 8).Value = 31449.611
$ws.Cells.Item(89, 9).Value = 57920.734
$ws.Cells.Item(89, 11).Value = 289603.67
$ws.Cells.Item(89, 13).Value = -283987.67

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(55, 8).Value = 35836100
$ws.Cells.Item(55, 10).Value = 1114566.6
$ws.Cells.Item(55, 12).Value = 3343699.8
$ws.Cells.Item(55, 14).Value = -3344053.8
$ws.Cells.Item(94, 8).Value = 13285
$ws.Cells.Item(94, 9).Value = 8997.5
$ws.Cells.Item(94, 10).Value = 15000
$ws.Cells.Item(94, 11).Value = 26992.5
$ws.Cells.Item(94, 12).Value = 45000
$ws.Cells.Item(94, 13).Value = -26316.5
$ws.Cells.Item(94, 14).Value = -46352
$ws.Cells.Item(96, 8).Value = 13833.167
$ws.Cells.Item(96, 10).Value = 13833.167
$ws.Cells.Item(96, 12).Value = 41499.501
$ws.Cells.Item(96, 14).Value = -45617.501
$ws.Cells.Item(100, 8).Value = 14661.667
$ws.Cells.Item(100, 9).Value = 0
$ws.Cells.Item(100, 10).Value = 14661.667
$ws.Cells.Item(100, 11).Value = 0
$ws.Cells.Item(100, 12).Value = 43985.001
$ws.Cells.Item(100, 13).ClearContents()
$ws.Cells.Item(100, 14).Value = -45607.001
$ws.Cells.Item(101, 8).Value = 6679652.5
$ws.Cells.Item(101, 10).Value = 6679652.5
$ws.Cells.Item(101, 12).Value = 20038957.5
$ws.Cells.Item(101, 14).Value = -20043825.5
$ws.Cells.Item(104, 8).Value = 1013
$ws.Cells.Item(104, 9).Value = 1013
$ws.Cells.Item(104, 11).Value = 3039
$ws.Cells.Item(104, 13).Value = -418
$ws.Cells.Item(109, 8).Value = 5211.85
$ws.Cells.Item(109, 9).Value = 1292.375
$ws.Cells.Item(109, 10).Value = 7824.8335
$ws.Cells.Item(109, 11).Value = 3877.125
$ws.Cells.Item(109, 12).Value = 23474.5005
$ws.Cells.Item(109, 13).Value = -2837.125
$ws.Cells.Item(109, 14).Value = -25554.5005
$ws.Cells.Item(112, 8).Value = 13556.25
$ws.Cells.Item(112, 10).Value = 14400
$ws.Cells.Item(112, 12).Value = 43200
$ws.Cells.Item(112, 14).Value = -45416
$ws.Cells.Item(115, 8).Value = 31892.5
$ws.Cells.Item(115, 9).Value = 31892.5
$ws.Cells.Item(115, 11).Value = 95677.5
$ws.Cells.Item(115, 13).Value = -94502.5
$ws.Cells.Item(118, 8).Value = 926.5
$ws.Cells.Item(118, 9).Value = 926.5
$ws.Cells.Item(118, 10).Value = 0
$ws.Cells.Item(118, 11).Value = 2779.5
$ws.Cells.Item(118, 12).Value = 0
$ws.Cells.Item(118, 13).Value = -1536.5
$ws.Cells.Item(118, 14).ClearContents()
$ws.Cells.Item(120, 8).Value = 15708.167
$ws.Cells.Item(120, 9).Value = 10499.714
$ws.Cells.Item(120, 10).Value = 23000
$ws.Cells.Item(120, 11).Value = 31499.142
$ws.Cells.Item(120, 12).Value = 69000
$ws.Cells.Item(120, 13).Value = -26661.142
$ws.Cells.Item(120, 14).Value = -78676

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(57, 8).Value = 20000
$ws.Cells.Item(57, 10).Value = 20000
$ws.Cells.Item(57, 12).Value = 20000
$ws.Cells.Item(57, 14).Value = -21640
$ws.Cells.Item(132, 8).Value = 9134.49
$ws.Cells.Item(132, 9).Value = 15069.875
$ws.Cells.Item(132, 11).Value = 45209.625
$ws.Cells.Item(132, 13).Value = -42679.625

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 3090.182
$ws.Cells.Item(46, 9).Value = 1437.375
$ws.Cells.Item(46, 10).Value = 4034.6428
$ws.Cells.Item(46, 11).Value = 1437.375
$ws.Cells.Item(46, 12).Value = 4034.6428
$ws.Cells.Item(46, 13).Value = -1249.375
$ws.Cells.Item(46, 14).Value = -4410.6428
$ws.Cells.Item(61, 8).Value = 13069.375
$ws.Cells.Item(61, 9).Value = 15166.667
$ws.Cells.Item(61, 10).Value = 6777.5
$ws.Cells.Item(61, 11).Value = 15166.667
$ws.Cells.Item(61, 12).Value = 6777.5
$ws.Cells.Item(61, 13).Value = -14964.667
$ws.Cells.Item(61, 14).Value = -7181.5
$ws.Cells.Item(93, 8).Value = 9498.4
$ws.Cells.Item(93, 9).Value = 7831
$ws.Cells.Item(93, 10).Value = 11999.5
$ws.Cells.Item(93, 11).Value = 7831
$ws.Cells.Item(93, 12).Value = 11999.5
$ws.Cells.Item(93, 13).Value = -6583
$ws.Cells.Item(93, 14).Value = -14495.5
$ws.Cells.Item(100, 8).Value = 3308
$ws.Cells.Item(100, 9).Value = 2646.8333
$ws.Cells.Item(100, 10).Value = 4299.75
$ws.Cells.Item(100, 11).Value = 2646.8333
$ws.Cells.Item(100, 12).Value = 4299.75
$ws.Cells.Item(100, 13).Value = -2105.8333
$ws.Cells.Item(100, 14).Value = -5381.75
$ws.Cells.Item(113, 8).Value = 13069.375
$ws.Cells.Item(113, 9).Value = 15166.667
$ws.Cells.Item(113, 10).Value = 6777.5
$ws.Cells.Item(113, 11).Value = 15166.667
$ws.Cells.Item(113, 12).Value = 6777.5
$ws.Cells.Item(113, 13).Value = -12996.667
$ws.Cells.Item(113, 14).Value = -11117.5
$ws.Cells.Item(132, 8).Value = 983637
$ws.Cells.Item(132, 9).Value = 1963659
$ws.Cells.Item(132, 10).Value = 3615
$ws.Cells.Item(132, 11).Value = 5890977
$ws.Cells.Item(132, 12).Value = 10845
$ws.Cells.Item(132, 13).Value = -5888447
$ws.Cells.Item(132, 14).Value = -15905

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 10561.6875
$ws.Cells.Item(62, 9).Value = 3061.25
$ws.Cells.Item(62, 11).Value = 3061.25
$ws.Cells.Item(62, 13).Value = -2437.25
$ws.Cells.Item(65, 8).Value = 10561.6875
$ws.Cells.Item(65, 9).Value = 3061.25
$ws.Cells.Item(65, 11).Value = 15306.25
$ws.Cells.Item(65, 13).Value = -12186.25
